$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (Strike#) values regenerated for column G, rows 2-22
$kValues = @(7, 8, 7, 7, 1, 3, 5, 3, 4, 9, 2, 5, 3, 7, 2, 3, 2, 7, 4, 2, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
